# Refresh the crypto price/volume table (rows 2-51) to the latest
# coinranking.com snapshot. Rows 17/18 also swap (TRON <-> WrappedEther
# changed rank order).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.309.62"
$ws.Range("E2").Value = "  -6.07%  "
$ws.Range("D3").Value = "3.049.27"
$ws.Range("E3").Value = "  -5.34%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'533.78"
$ws.Range("E5").Value = "  -7.55%  "
$ws.Range("D6").Value = "'130.31"
$ws.Range("E6").Value = "  -14.21%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "3.043.50"
$ws.Range("E8").Value = "  -5.37%  "
$ws.Range("E9").Value = "  -5.46%  "
$ws.Range("D10").Value = "'0.151"
$ws.Range("E10").Value = "  -6.73%  "
$ws.Range("E11").Value = "  -13.33%  "
$ws.Range("D12").Value = "'0.452"
$ws.Range("E12").Value = "  -7.21%  "
$ws.Range("D13").Value = "'0.0000221"
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("D14").Value = "'33.75"
$ws.Range("E14").Value = "  -12.10%  "
$ws.Range("D15").Value = "3.518.32"
$ws.Range("E15").Value = "  -5.38%  "
$ws.Range("D16").Value = "62.441.32"
$ws.Range("E16").Value = "  -5.80%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.110"
$ws.Range("E17").Value = "  -3.90%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.069.90"
$ws.Range("E18").Value = "  -4.70%  "
$ws.Range("D19").Value = "'6.50"
$ws.Range("E19").Value = "  -8.91%  "
$ws.Range("D20").Value = "'473.89"
$ws.Range("E20").Value = "  -11.79%  "
$ws.Range("D21").Value = "'13.10"
$ws.Range("E21").Value = "  -10.22%  "
$ws.Range("D22").Value = "'0.689"
$ws.Range("E22").Value = "  -7.27%  "
$ws.Range("D23").Value = "'7.06"
$ws.Range("E23").Value = "  -8.89%  "
$ws.Range("D24").Value = "'77.78"
$ws.Range("E24").Value = "  -4.23%  "
$ws.Range("D25").Value = "'11.80"
$ws.Range("E25").Value = "  -12.84%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -9.79%  "
$ws.Range("D28").Value = "'8.03"
$ws.Range("E28").Value = "  -14.76%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "'25.39"
$ws.Range("E30").Value = "  -7.78%  "
$ws.Range("D31").Value = "'1.85"
$ws.Range("E31").Value = "  -18.00%  "
$ws.Range("D32").Value = "'1.08"
$ws.Range("E32").Value = "  -7.97%  "
$ws.Range("D33").Value = "'58.17"
$ws.Range("E33").Value = "  +5.84%  "
$ws.Range("D34").Value = "'2.36"
$ws.Range("E34").Value = "  -14.32%  "
$ws.Range("D35").Value = "'5.84"
$ws.Range("E35").Value = "  -8.33%  "
$ws.Range("D36").Value = "'5.14"
$ws.Range("E36").Value = "  -8.63%  "
$ws.Range("D37").Value = "'462.48"
$ws.Range("E37").Value = "  -17.72%  "
$ws.Range("D38").Value = "3.086.41"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("D39").Value = "'0.0384"
$ws.Range("E39").Value = "  -15.04%  "
$ws.Range("D40").Value = "'0.0777"
$ws.Range("E40").Value = "  -9.63%  "
$ws.Range("D41").Value = "'0.112"
$ws.Range("E41").Value = "  -13.01%  "
$ws.Range("D42").Value = "'7.94"
$ws.Range("E42").Value = "  -7.77%  "
$ws.Range("E43").Value = "  -15.29%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'0.245"
$ws.Range("E45").Value = "  -13.78%  "
$ws.Range("E46").Value = "  -16.23%  "
$ws.Range("D47").Value = "'24.02"
$ws.Range("E47").Value = "  -9.88%  "
$ws.Range("D48").Value = "'116.38"
$ws.Range("E48").Value = "  -6.67%  "
$ws.Range("E49").Value = "  -6.03%  "
$ws.Range("D50").Value = "0.0₃0503"
$ws.Range("E50").Value = "  -9.43%  "
$ws.Range("D51").Value = "'1.95"
$ws.Range("E51").Value = "  -11.91%  "
